# Auto-generated Excel COM-interop script applying the Famfrit_Profits.xlsx diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("N25").ClearContents()

$ws.Range("H32").Value = 5275.3335
$ws.Range("I32").Value = 3624.75
$ws.Range("J32").Value = 6595.8
$ws.Range("K32").Value = 3624.75
$ws.Range("L32").Value = 6595.8
$ws.Range("M32").Value = -3298.75
$ws.Range("N32").Value = -7247.8

$ws.Range("H40").Value = 2819.4
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 3499
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 3499
$ws.Range("M40").Value = -1625

$ws.Range("H41").Value = 272.36365
$ws.Range("I41").Value = 399
$ws.Range("J41").Value = 259.7
$ws.Range("K41").Value = 399
$ws.Range("L41").Value = 259.7
$ws.Range("M41").Value = 41
$ws.Range("N41").Value = -1139.7

$ws.Range("H43").Value = 2503149.5
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 2503149.5
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 2503149.5
$ws.Range("N43").Value = -2503287.5

$ws.Range("H57").Value = 150000
$ws.Range("I57").Value = 50001
$ws.Range("J57").Value = 199999.5
$ws.Range("K57").Value = 150003
$ws.Range("L57").Value = 599998.5
$ws.Range("M57").Value = -149504
$ws.Range("N57").Value = -600996.5

$ws.Range("H112").Value = 5265962.5
$ws.Range("I112").Value = 2323
$ws.Range("J112").Value = 5558387
$ws.Range("K112").Value = 6969
$ws.Range("L112").Value = 16675161
$ws.Range("M112").Value = -5861

$ws.Range("H113").Value = 5626.5625
$ws.Range("I113").Value = 4693.4165
$ws.Range("J113").Value = 8426
$ws.Range("K113").Value = 4693.4165
$ws.Range("L113").Value = 8426
$ws.Range("M113").Value = -1439.4165

$ws.Range("H116").Value = 7476.706
$ws.Range("I116").Value = 7513.7144
$ws.Range("J116").Value = 7450.8
$ws.Range("K116").Value = 7513.7144
$ws.Range("L116").Value = 7450.8
$ws.Range("M116").Value = -4071.7144

$ws.Range("H132").Value = 1238.2
$ws.Range("I132").Value = 1059.4082
$ws.Range("J132").Value = 9999
$ws.Range("K132").Value = 3178.2246
$ws.Range("L132").Value = 29997
$ws.Range("M132").Value = -648.2246000000005
$ws.Range("N132").Value = -35057

$ws.Range("H137").Value = 44448092
$ws.Range("I137").Value = 38465016
$ws.Range("J137").Value = 52635456
$ws.Range("K137").Value = 115395048
$ws.Range("L137").Value = 157906368
$ws.Range("M137").Value = -115392498
$ws.Range("N137").Value = -157911468

$ws.Range("H138").Value = 6543251
$ws.Range("I138").Value = 3422.1765
$ws.Range("J138").Value = 9813165
$ws.Range("K138").Value = 10266.5295
$ws.Range("L138").Value = 29439495
$ws.Range("M138").Value = -5126.529500000001
$ws.Range("N138").Value = -29449775

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 837.63635
$ws.Range("I2").Value = 837.63635
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 837.63635
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -724.63635
$ws.Range("N2").ClearContents()

$ws.Range("H5").Value = 218.33333
$ws.Range("I5").Value = 144.66667
$ws.Range("J5").Value = 365.66666
$ws.Range("K5").Value = 144.66667
$ws.Range("L5").Value = 365.66666
$ws.Range("M5").Value = -32.66667000000001
$ws.Range("N5").Value = -589.66666

$ws.Range("H8").Value = 334340000
$ws.Range("I8").Value = 1000000000
$ws.Range("J8").Value = 1509999.5
$ws.Range("K8").Value = 1000000000
$ws.Range("L8").Value = 1509999.5
$ws.Range("M8").Value = -999999856
$ws.Range("N8").Value = -1510287.5

$ws.Range("H32").Value = 14713368
$ws.Range("I32").Value = 16134452
$ws.Range("J32").Value = 28829
$ws.Range("K32").Value = 16134452
$ws.Range("L32").Value = 28829
$ws.Range("M32").Value = -16134165
$ws.Range("N32").Value = -29403

$ws.Range("H61").Value = 41671400
$ws.Range("I61").Value = 55558452
$ws.Range("J61").Value = 10250
$ws.Range("K61").Value = 55558452
$ws.Range("L61").Value = 10250
$ws.Range("M61").Value = -55558240

$ws.Range("H74").Value = 166855180
$ws.Range("I74").Value = 250281780
$ws.Range("J74").Value = 1974.5
$ws.Range("K74").Value = 250281780
$ws.Range("L74").Value = 1974.5
$ws.Range("M74").Value = -250280906
$ws.Range("N74").Value = -3722.5

$ws.Range("H77").Value = 166855180
$ws.Range("I77").Value = 250281780
$ws.Range("J77").Value = 1974.5
$ws.Range("K77").Value = 1251408900
$ws.Range("L77").Value = 9872.5
$ws.Range("M77").Value = -1251404532
$ws.Range("N77").Value = -18608.5

$ws.Range("H102").Value = 3503.8572
$ws.Range("I102").Value = 3532.3333
$ws.Range("J102").Value = 3333
$ws.Range("K102").Value = 3532.3333
$ws.Range("L102").Value = 3333
$ws.Range("M102").Value = -1910.3333
$ws.Range("N102").Value = -6577

$ws.Range("H116").Value = 837.63635
$ws.Range("I116").Value = 837.63635
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 837.63635
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1456.36365
$ws.Range("N116").ClearContents()

$ws.Range("H132").Value = 34486228
$ws.Range("I132").Value = 3920.0952
$ws.Range("J132").Value = 125002280
$ws.Range("K132").Value = 11760.2856
$ws.Range("L132").Value = 375006840
$ws.Range("M132").Value = -9230.285600000001

$ws.Range("H134").Value = 444999.5
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 444999.5
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 444999.5
$ws.Range("N134").Value = -455139.5

$ws.Range("H136").Value = 41671400
$ws.Range("I136").Value = 55558452
$ws.Range("J136").Value = 10250
$ws.Range("K136").Value = 166675356
$ws.Range("L136").Value = 30750
$ws.Range("M136").Value = -166672806

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 837.63635
$ws.Range("I3").Value = 837.63635
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 837.63635
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -723.63635
$ws.Range("N3").ClearContents()

$ws.Range("H4").Value = 218.33333
$ws.Range("I4").Value = 144.66667
$ws.Range("J4").Value = 365.66666
$ws.Range("K4").Value = 144.66667
$ws.Range("L4").Value = 365.66666
$ws.Range("M4").Value = -29.66667000000001
$ws.Range("N4").Value = -595.66666

$ws.Range("H105").Value = 1530.8422
$ws.Range("I105").Value = 1477.6428
$ws.Range("J105").Value = 1679.8
$ws.Range("K105").Value = 1477.6428
$ws.Range("L105").Value = 1679.8
$ws.Range("M105").Value = 269.3571999999999
$ws.Range("N105").Value = -5173.8

$ws.Range("H134").Value = 1580.9828
$ws.Range("I134").Value = 1569.1569
$ws.Range("J134").Value = 1667.1428
$ws.Range("K134").Value = 4707.4707
$ws.Range("L134").Value = 5001.428400000001
$ws.Range("M134").Value = -2172.4707

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 3015
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 3015
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 3015
$ws.Range("N21").Value = -3485
$ws.Range("M21").ClearContents()

$ws.Range("H26").Value = 500
$ws.Range("I26").Value = 0
$ws.Range("J26").Value = 500
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = 500
$ws.Range("N26").Value = -1074

$ws.Range("H29").Value = 2000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 2000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 2000
$ws.Range("N29").Value = -2586

$ws.Range("H31").Value = 38467388
$ws.Range("I31").Value = 4398.933
$ws.Range("J31").Value = 90916920
$ws.Range("K31").Value = 4398.933
$ws.Range("L31").Value = 90916920
$ws.Range("M31").Value = -4103.933
$ws.Range("N31").Value = -90917510

$ws.Range("H34").Value = 38467388
$ws.Range("I34").Value = 4398.933
$ws.Range("J34").Value = 90916920
$ws.Range("K34").Value = 4398.933
$ws.Range("L34").Value = 90916920
$ws.Range("M34").Value = -4196.933
$ws.Range("N34").Value = -90917324

$ws.Range("H132").Value = 3279.261
$ws.Range("I132").Value = 3203.6316
$ws.Range("J132").Value = 3638.5
$ws.Range("K132").Value = 9610.8948
$ws.Range("L132").Value = 10915.5
$ws.Range("M132").Value = -7080.8948
$ws.Range("N132").Value = -15975.5

$ws.Range("H134").Value = 1237.9166
$ws.Range("I134").Value = 1142.1875
$ws.Range("J134").Value = 2003.75
$ws.Range("K134").Value = 3426.5625
$ws.Range("L134").Value = 6011.25
$ws.Range("M134").Value = -891.5625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H32").Value = 1887.8572
$ws.Range("I32").Value = 360
$ws.Range("J32").Value = 2499
$ws.Range("K32").Value = 1080
$ws.Range("L32").Value = 7497
$ws.Range("M32").Value = -797
$ws.Range("N32").Value = -8063

$ws.Range("H113").Value = 3728.842
$ws.Range("I113").Value = 2699.75
$ws.Range("J113").Value = 4003.2666
$ws.Range("K113").Value = 8099.25
$ws.Range("L113").Value = 12009.7998
$ws.Range("M113").Value = -5929.25

$ws.Range("H132").Value = 2302053
$ws.Range("I132").Value = 1714.8334
$ws.Range("J132").Value = 2902141.2
$ws.Range("K132").Value = 15433.5006
$ws.Range("L132").Value = 26119270.8
$ws.Range("M132").Value = -12903.5006

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 998.2222
$ws.Range("I107").Value = 1010.5
$ws.Range("J107").Value = 900
$ws.Range("K107").Value = 1010.5
$ws.Range("L107").Value = 900
$ws.Range("M107").Value = 909.5
$ws.Range("N107").Value = -4740

$ws.Range("H126").Value = 5268806
$ws.Range("I126").Value = 2783371.2
$ws.Range("J126").Value = 9529551
$ws.Range("K126").Value = 8350113.600000001
$ws.Range("L126").Value = 28588653
$ws.Range("M126").Value = -8347643.600000001
$ws.Range("N126").Value = -28593593

$ws.Range("H132").Value = 3795.8298
$ws.Range("I132").Value = 2674.2942
$ws.Range("J132").Value = 6729.077
$ws.Range("K132").Value = 8022.882599999999
$ws.Range("L132").Value = 20187.231
$ws.Range("M132").Value = -5492.882599999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4666.9565
$ws.Range("I7").Value = 4393.357
$ws.Range("J7").Value = 5092.5557
$ws.Range("K7").Value = 4393.357
$ws.Range("L7").Value = 5092.5557
$ws.Range("M7").Value = -4281.357

$ws.Range("H16").Value = 3001
$ws.Range("I16").Value = 3001
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 3001
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2831

$ws.Range("H68").Value = 5499.8
$ws.Range("I68").Value = 3333.3333
$ws.Range("J68").Value = 6428.2856
$ws.Range("K68").Value = 3333.3333
$ws.Range("L68").Value = 6428.2856
$ws.Range("M68").Value = -2584.3333
$ws.Range("N68").Value = -7926.2856

$ws.Range("H71").Value = 5499.8
$ws.Range("I71").Value = 3333.3333
$ws.Range("J71").Value = 6428.2856
$ws.Range("K71").Value = 16666.6665
$ws.Range("L71").Value = 32141.428
$ws.Range("M71").Value = -12922.6665
$ws.Range("N71").Value = -39629.428

$ws.Range("H100").Value = 3269.3667
$ws.Range("I100").Value = 2905.6
$ws.Range("J100").Value = 3633.1333
$ws.Range("K100").Value = 2905.6
$ws.Range("L100").Value = 3633.1333
$ws.Range("M100").Value = -2364.6

$ws.Range("H126").Value = 4666.9565
$ws.Range("I126").Value = 4393.357
$ws.Range("J126").Value = 5092.5557
$ws.Range("K126").Value = 13180.071
$ws.Range("L126").Value = 15277.6671
$ws.Range("M126").Value = -10710.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()

$ws.Range("H126").Value = 6521.1577
$ws.Range("I126").Value = 7021.643
$ws.Range("J126").Value = 5119.8
$ws.Range("K126").Value = 21064.929
$ws.Range("L126").Value = 15359.4
$ws.Range("M126").Value = -18594.929
$ws.Range("N126").Value = -20299.4

$ws.Range("H132").Value = 4027.25
$ws.Range("I132").Value = 4003.2
$ws.Range("J132").Value = 4129
$ws.Range("K132").Value = 12009.6
$ws.Range("L132").Value = 12387
$ws.Range("M132").Value = -9479.599999999999
$ws.Range("N132").Value = -17447
